$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# V14 2014-Høst (row 14): mark 1e as also not curriculum -> "1f" becomes "1e-f"
$ws.Range("D14").Value = "1e-f"

# V15 2015-Vår (row 15): mark 1c as not curriculum (previously empty)
$ws.Range("D15").Value = "1c"

# Update selection to reflect last-edited cell
$ws.Range("D15").Select()
